$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the header labels to lowercase / ASCII-only variants
$ws.Range("A1").Value = "jahr"
$ws.Range("B1").Value = "fluege"
$ws.Range("C1").Value = "fluggaeste"
$ws.Range("D1").Value = "sitzladefaktor"

# Clear the explicit number format style that was applied on C9:C13 so
# they fall back to the default "General" style again.
$ws.Range("C9:C13").ClearFormats()

# Move the active selection from F12 to G10
$ws.Range("G10").Select()
